# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Re-rank three pairs of countries whose case counts crossed over
#   (Australia/Chequia, Panama/Republica Dominicana, Estonia/Irak)
# - Refresh numeric COVID-19 stats for a number of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 18:22"

# --- Country re-ranking (swap display names, rank order in the table) -
$ws.Range("A37").Value = "Chequia"
$ws.Range("A38").Value = "Australia"

$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("A49").Value = "Panama"

$ws.Range("A68").Value = "Irak"
$ws.Range("A69").Value = "Estonia"

# --- Row 4 (Estados Unidos) --------------------------------------------
$ws.Range("B4").Value = 715105
$ws.Range("C4").Value = 5370
$ws.Range("D4").Value = 63841
$ws.Range("E4").Value = 613375
$ws.Range("G4").Value = 735
$ws.Range("H4").Value = 37889

# --- Row 6 (Italia) -----------------------------------------------------
$ws.Range("B6").Value = 175925
$ws.Range("C6").Value = 3491
$ws.Range("D6").Value = 44927
$ws.Range("E6").Value = 107771
$ws.Range("F6").Value = 2733
$ws.Range("G6").Value = 482
$ws.Range("H6").Value = 23227

# --- Row 8 (Alemania) ----------------------------------------------------
$ws.Range("B8").Value = 142569
$ws.Range("C8").Value = 1172
$ws.Range("E8").Value = 52764
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = 4405

# --- Row 15 (Brasil) -----------------------------------------------------
$ws.Range("B15").Value = 34913
$ws.Range("C15").Value = 1231
$ws.Range("E15").Value = 18686
$ws.Range("G15").Value = 60
$ws.Range("H15").Value = 2201

# --- Row 37 (now Chequia) -------------------------------------------------
$ws.Range("B37").Value = 6606
$ws.Range("C37").Value = 57
$ws.Range("D37").Value = 1227
$ws.Range("E37").Value = 5198
$ws.Range("F37").Value = 86
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 181

# --- Row 38 (now Australia) -----------------------------------------------
$ws.Range("B38").Value = 6565
$ws.Range("C38").Value = 32
$ws.Range("D38").Value = 4163
$ws.Range("E38").Value = 2333
$ws.Range("F38").Value = 55
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 69

# --- Row 48 (now Republica Dominicana) ------------------------------------
$ws.Range("B48").Value = 4335
$ws.Range("C48").Value = 209
$ws.Range("D48").Value = 312
$ws.Range("E48").Value = 3806
$ws.Range("F48").Value = 146
$ws.Range("G48").Value = 17
$ws.Range("H48").Value = 217

# --- Row 49 (now Panama) --------------------------------------------------
$ws.Range("B49").Value = 4210
$ws.Range("C49").Value = 194
$ws.Range("D49").Value = 122
$ws.Range("E49").Value = 3972
$ws.Range("F49").Value = 96
$ws.Range("G49").Value = 7
$ws.Range("H49").Value = 116

# --- Row 58 (Argelia) ----------------------------------------------------
$ws.Range("B58").Value = 2534
$ws.Range("C58").Value = 116
$ws.Range("D58").Value = 894
$ws.Range("E58").Value = 1273
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 367

# --- Row 68 (now Irak) -----------------------------------------------------
$ws.Range("B68").Value = 1513
$ws.Range("C68").Value = 31
$ws.Range("D68").Value = 953
$ws.Range("E68").Value = 478
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 82

# --- Row 69 (now Estonia) ---------------------------------------------------
$ws.Range("B69").Value = 1512
$ws.Range("C69").Value = 53
$ws.Range("D69").Value = 162
$ws.Range("E69").Value = 1312
$ws.Range("F69").Value = 11
$ws.Range("H69").Value = 38

# --- Row 82 (Cuba) ------------------------------------------------------
$ws.Range("F82").Value = 15

# --- Row 93 (Libano) ------------------------------------------------------
$ws.Range("D93").Value = 99
$ws.Range("E93").Value = 552
